# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# Before:  Sheets = "2021-Q4", "总计"
# After:   Sheets = "2021-Q4", "2022-Q1", "总计"
#   - "2022-Q1" is a brand-new fund-holdings sheet (same layout as "2021-Q4")
#   - "总计" gains a new summary row for 2022-Q1 (inserted above 2021-Q4's row)
#
# NOTE: worksheet variables captured before a sheet is inserted/removed can
# resolve to the wrong tab afterwards (collection position shifts under
# them), so "总计" is always re-fetched by name right before it is used.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# Remember the old totals-sheet values before we touch anything (.Value2 is
# used for reads -- .Value's getter is unreliable in this host).
$total = $wb.Worksheets.Item("总计")
$oldDateLabel = $total.Range("B2").Value2
$oldCount = $total.Range("C2").Value2
$oldValue = $total.Range("D2").Value2

# Drop the old "总计" sheet now (before any new sheet is created) so that the
# sheet re-added in step 2 below is allocated the next sheetId in sequence,
# matching 1/2/3 for 2021-Q4/2022-Q1/总计 like a from-scratch save would.
$total.Delete()

# ---------------------------------------------------------------------------
# 1. Build the "2022-Q1" sheet by cloning "2021-Q4" (so headers/column
#    widths/fonts/borders all come along for free) and dropping the clone
#    right after "2021-Q4".
# ---------------------------------------------------------------------------
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

# New fund-holdings rows for 2022-Q1 (code, name, fund size, stock position,
# position ratio, held market value, position rank).
$rows = @(
    @("501011", "汇添富中证中药指数（LOF）A", "24.13", "94.57", "2.93", "0.7070", 10),
    @("710001", "富安达优势成长混合", "9.04", "94.45", "4.57", "0.4131", 5),
    @("501012", "汇添富中证中药指数（LOF）C", "8.91", "94.57", "2.93", "0.2611", 10),
    @("011269", "中银证券优势制造股票型证券投资基金A", "1.39", "93.51", "4.25", "0.0591", 9),
    @("010746", "富安达长三角区域主题混合", "1.16", "92.45", "4.79", "0.0556", 7),
    @("159804", "国寿安保国证创业板中盘精选88ETF", "2.10", "98.79", "2.43", "0.0510", 3),
    @("005293", "诺德新旺灵活配置混合", "0.44", "89.98", "5.39", "0.0237", 5),
    @("009789", "富安达科技创新混合", "0.56", "94.07", "3.80", "0.0213", 9),
    @("540007", "汇丰晋信中小盘股票", "0.61", "93.28", "2.07", "0.0126", 8),
    @("011270", "中银证券优势制造股票型证券投资基金C", "0.21", "93.51", "4.25", "0.0089", 9)
)

# "2021-Q4" only had 8 data rows (rows 2-9); 2022-Q1 needs 10 (rows 2-11), so
# stamp the bold/bordered index-column style onto the two extra rows before
# writing into them.
$q1.Range("A9").Copy() | Out-Null
$q1.Range("A10:A11").PasteSpecial(-4122)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q1.Range("A$r").Value2 = $i

    # Columns B:G are stored as literal text (fund codes/ratios keep their
    # printed form, e.g. leading zeros / trailing zeros), so force text
    # format before assigning, then strip the format stamp back off so the
    # cells end up unstyled just like the source sheet.
    $q1.Range("B$r`:G$r").NumberFormat = "@"
    $q1.Range("B$r").Value2 = $row[0]
    $q1.Range("C$r").Value2 = $row[1]
    $q1.Range("D$r").Value2 = $row[2]
    $q1.Range("E$r").Value2 = $row[3]
    $q1.Range("F$r").Value2 = $row[4]
    $q1.Range("G$r").Value2 = $row[5]
    $q1.Range("B$r`:G$r").ClearFormats()

    $q1.Range("H$r").Value2 = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Rebuild "总计" so it picks up sheetId/rId 3 (i.e. it becomes the last
#    sheet again) while keeping its visual formatting, then fill it with the
#    2022-Q1 row on top of the carried-over 2021-Q4 row.
# ---------------------------------------------------------------------------
$total2 = $wb.Worksheets.Add($null, $q1)
$total2.Name = "总计"

# Re-apply the bold/bordered style used for the header row + index column,
# copying it straight from the 2022-Q1 sheet we just finished building.
$q1.Range("B1").Copy() | Out-Null
$total2.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy() | Out-Null
$total2.Range("A2:A3").PasteSpecial(-4122)

$total2.Range("B1").Value2 = "日期"
$total2.Range("C1").Value2 = "持有数量(只)"
$total2.Range("D1").Value2 = "持有市值(亿元)"

$total2.Range("A2").Value2 = 0
$total2.Range("B2").Value2 = "2022-Q1"
$total2.Range("C2").Value2 = 10
$total2.Range("D2").Value2 = 1.61

$total2.Range("A3").Value2 = 1
$total2.Range("B3").Value2 = $oldDateLabel
$total2.Range("C3").Value2 = $oldCount
$total2.Range("D3").Value2 = $oldValue
